# Pathways Update - add animal records to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: Name --- (entered top-to-bottom; A3 intentionally left blank)
$ws.Range("A2").Value = "Wishbone"
$ws.Range("A4").Value = "Benji"
$ws.Range("A5").Value = "Magnus"
$ws.Range("A6").Value = "Scruffles"
$ws.Range("A7").Value = "Huey"
$ws.Range("A8").Value = "PUFF"

# --- Column B: AID (kept as text so leading context/zeros aren't lost) ---
$ws.Range("B2:B8").NumberFormat = "@"
$ws.Range("B2").Value = "57230680"
$ws.Range("B3").Value = "58622374"
$ws.Range("B4").Value = "58639325"
$ws.Range("B5").Value = "58690617"
$ws.Range("B6").Value = "58419285"
$ws.Range("B7").Value = "58903460"
$ws.Range("B8").Value = "58940639"
$ws.Range("B2:B8").Style = "Normal"

# --- Column C: Species ---
$ws.Range("C2").Value = "Dog"
$ws.Range("C3").Value = "Rabbit"
$ws.Range("C4").Value = "Dog"
$ws.Range("C5").Value = "Dog"
$ws.Range("C6").Value = "Dog"
$ws.Range("C7").Value = "Dog"
$ws.Range("C8").Value = "Reptile/Amphibian"

# --- Column D: Location ---
$ws.Range("D2").Value = "Dog E"
$ws.Range("D3").Value = "Foster Home"
$ws.Range("D4").Value = "Dog Adoptions A"
$ws.Range("D5").Value = "If The Fur Fits"
$ws.Range("D6").Value = "Foster Home"
$ws.Range("D7").Value = "Dog Adoptions B"
$ws.Range("D8").Value = "Foster Home"

# --- Column E: Intake Date (kept as text to preserve the mixed date formats) ---
$ws.Range("E2:E8").NumberFormat = "@"
$ws.Range("E2").Value = "11/13/2024"
$ws.Range("E3").Value = "6/2/2025"
$ws.Range("E4").Value = "6/4/2025"
$ws.Range("E5").Value = "6/19/2025"
$ws.Range("E6").Value = "5/1/2025"
$ws.Range("E7").Value = "7/17/25"
$ws.Range("E8").Value = "7/19/25"
$ws.Range("E2:E8").Style = "Normal"

# --- Resize columns to fit the new content (matches the widths Excel's
#     "AutoFit Column Width" produced for this data) ---
$ws.Columns("A:E").AutoFit()
$ws.Columns("A").ColumnWidth = 9.166666666666666
$ws.Columns("B").ColumnWidth = 8.166666666666666
$ws.Columns("C").ColumnWidth = 17.5
$ws.Columns("D").ColumnWidth = 15.0
$ws.Columns("E").ColumnWidth = 10.333333333333332

# --- Restore the selection left by the author ---
$ws.Range("C7").Select() | Out-Null
